# TestSuite.xlsx edit: add new "CreateLead01" test steps (Lead status,
# Name salutation, Company name, Save button) and update existing step
# data (login creds, button labels, run/wait flags), per commit:
#   "Addition on TestSuite excel file - New lead creation test cases
#    added - Objectrepository updated with xpaths - selectItem() method
#    added to Keywords class"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestLeadSteps")
$ws.Activate() | Out-Null

# --- Prepare formatting for the 4 brand-new rows (12-15) by cloning the
# formats (fill/border/numberformat) of an existing data row (columns
# A-H only), so the new rows visually match the rest of the table
# (incl. column E's distinct style).
$ws.Range("A9:H9").Copy() | Out-Null
$ws.Range("A12:H12").PasteSpecial(-4122) | Out-Null
$ws.Range("A9:H9").Copy() | Out-Null
$ws.Range("A13:H13").PasteSpecial(-4122) | Out-Null
$ws.Range("A9:H9").Copy() | Out-Null
$ws.Range("A14:H14").PasteSpecial(-4122) | Out-Null
$ws.Range("A9:H9").Copy() | Out-Null
$ws.Range("A15:H15").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Row 4: "Enter username" now uses the real login email; it also
# becomes a mailto hyperlink, keeping the same wait flag as row 3 below.
$ws.Range("H4").Value = $false
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:lavanya.kpv@gmail.com", "", "", "lavanya.kpv@gmail.com") | Out-Null

# --- Row 5: "Enter password" uses the real login password.
$ws.Range("D5").Value = "kavilputhenlavu1"
$ws.Range("H5").Value = $false

# --- Rows 6-9: waitmode flag now off for all of them; row 8's runmode
# turned on, and row 9's description typo is fixed with runmode turned on.
$ws.Range("H6").Value = $false
$ws.Range("H7").Value = $false
$ws.Range("G8").Value = $true
$ws.Range("H8").Value = $false
$ws.Range("B9").Value = "Click on new button"
$ws.Range("G9").Value = $true
$ws.Range("H9").Value = $false

# --- Rows 10-11: repurposed from the old "Enter First/Last Name" rows
# into two brand-new lead-creation steps (select lead status / salutation).
$ws.Range("B10").Value = "Select Lead status"
$ws.Range("C10").Value = "selectItem"
$ws.Range("D10").Value = "New"
$ws.Range("F10").Value = "leadstatselector"
$ws.Range("G10").Value = $true
$ws.Range("H10").Value = $false

$ws.Range("B11").Value = "Select Name salutation"
$ws.Range("C11").Value = "selectItem"
$ws.Range("D11").Value = "Mr."
$ws.Range("F11").Value = "namesalutselector"
$ws.Range("G11").Value = $true
$ws.Range("H11").Value = $false

# --- Rows 12-15: new rows continuing the lead-creation flow (first name,
# last name, company name, save button).
$ws.Range("A12").Value = "CreateLead01"
$ws.Range("B12").Value = "Enter First Name"
$ws.Range("C12").Value = "inputData"
$ws.Range("D12").Value = "Larry"
$ws.Range("E12").Value = "xpath"
$ws.Range("F12").Value = "firstnametextbox"
$ws.Range("G12").Value = $true
$ws.Range("H12").Value = $false

$ws.Range("A13").Value = "CreateLead01"
$ws.Range("B13").Value = "Enter Last Name"
$ws.Range("C13").Value = "inputData"
$ws.Range("D13").Value = "Page"
$ws.Range("E13").Value = "xpath"
$ws.Range("F13").Value = "lastnametextbox"
$ws.Range("G13").Value = $true
$ws.Range("H13").Value = $false

$ws.Range("A14").Value = "CreateLead01"
$ws.Range("B14").Value = "Enter Company Name"
$ws.Range("C14").Value = "inputData"
$ws.Range("D14").Value = "Sutherland"
$ws.Range("E14").Value = "xpath"
$ws.Range("F14").Value = "companytextbox"
$ws.Range("G14").Value = $true
$ws.Range("H14").Value = $false

$ws.Range("A15").Value = "CreateLead01"
$ws.Range("B15").Value = "Click on Save button"
$ws.Range("C15").Value = "clickElement"
$ws.Range("D15").Value = "null"
$ws.Range("E15").Value = "xpath"
$ws.Range("F15").Value = "savebutton"
$ws.Range("G15").Value = $true
$ws.Range("H15").Value = $false

# --- Column F needs to be a bit wider to fit the new xpath selector names.
$ws.Columns.Item(6).ColumnWidth = 18.14

# --- Final UI state: active cell moves to D5 (matches the edited sheet).
$ws.Range("D5").Select() | Out-Null
